$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "96.850.15"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  +0.63%  "
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.686.71"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  +1.00%  "
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  -0.05%  "
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "236.26"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  -1.83%  "
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +3.37%  "
$r.Style = "Normal"

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "654.53"
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -1.06%  "
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.423"
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  +0.46%  "
$r.Style = "Normal"

$r = $ws.Range("B9")
$r.NumberFormat = "@"
$r.Value = "Cardano"
$r.Style = "Normal"

$r = $ws.Range("C9")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "1.07"
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  -1.30%  "
$r.Style = "Normal"

$r = $ws.Range("B10")
$r.NumberFormat = "@"
$r.Value = "USDC"
$r.Style = "Normal"

$r = $ws.Range("C10")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  -0.04%  "
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "3.684.95"
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  +1.12%  "
$r.Style = "Normal"

$r = $ws.Range("B12")
$r.NumberFormat = "@"
$r.Value = "Avalanche"
$r.Style = "Normal"

$r = $ws.Range("C12")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "44.03"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -1.25%  "
$r.Style = "Normal"

$r = $ws.Range("B13")
$r.NumberFormat = "@"
$r.Value = "TRON"
$r.Style = "Normal"

$r = $ws.Range("C13")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.208"
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  +2.18%  "
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.0000298"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  +11.27%  "
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "6.73"
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +2.06%  "
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "4.372.81"
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +1.01%  "
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "96.663.03"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  +0.52%  "
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "9.00"
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +2.28%  "
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "3.683.69"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  +1.01%  "
$r.Style = "Normal"

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "12.96"
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +2.61%  "
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "18.59"
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +2.32%  "
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.509"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -4.34%  "
$r.Style = "Normal"

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "520.84"
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  +0.10%  "
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "3.40"
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -0.78%  "
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  +3.47%  "
$r.Style = "Normal"

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "6.88"
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  +0.68%  "
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.204"
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  +23.49%  "
$r.Style = "Normal"

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "101.00"
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  -1.05%  "
$r.Style = "Normal"

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "13.32"
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  +3.33%  "
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "12.36"
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  +2.49%  "
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  -0.76%  "
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.998"
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  +0.41%  "
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  +2.45%  "
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "1.84"
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  +1.56%  "
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.998"
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "32.12"
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  -1.22%  "
$r.Style = "Normal"

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "642.49"
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  +3.99%  "
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  +1.55%  "
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "8.78"
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +1.25%  "
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +0.01%  "
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "6.82"
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  +11.56%  "
$r.Style = "Normal"

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "40.83"
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  -4.76%  "
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.04"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  +5.89%  "
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  +1.25%  "
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.949"
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +0.47%  "
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.457"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  +7.46%  "
$r.Style = "Normal"

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0456"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  +1.21%  "
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  +0.13%  "
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  -0.46%  "
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "8.51"
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +0.72%  "
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -0.44%  "
$r.Style = "Normal"

